$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add D1=3, E1=4, copying style from existing header cells ---
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Updated values for column C (rows 2-12) ---
$ws.Range("C2").Value = -5.04614857794682
$ws.Range("C3").Value = -1.194025718115943
$ws.Range("C4").Value = -0.07161795042852842
$ws.Range("C5").Value = -0.4111325302719243
$ws.Range("C6").Value = 0.0159267162195228
$ws.Range("C7").Value = 0.1002874912444511
$ws.Range("C8").Value = 0.1290666877551792
$ws.Range("C9").Value = 0.0276633633304105
$ws.Range("C10").Value = 0.03180697780879011
$ws.Range("C11").Value = 0.005376147938177376
$ws.Range("C12").Value = 0.03799774138790459

# --- New values for column D (rows 2-12) ---
$ws.Range("D2").Value = -4.703688665936504
$ws.Range("D3").Value = -1.184465726122439
$ws.Range("D4").Value = 0.02172012643327927
$ws.Range("D5").Value = -0.1546895415091126
$ws.Range("D6").Value = -0.04401287182054063
$ws.Range("D7").Value = 0.06273304633550032
$ws.Range("D8").Value = 0.03477873499536071
$ws.Range("D9").Value = 0.0247399727427734
$ws.Range("D10").Value = 0.02298403018683341
$ws.Range("D11").Value = 0.008014114575528983
$ws.Range("D12").Value = 0.02684937632215087

# --- New values for column E (rows 2-12) ---
$ws.Range("E2").Value = -4.331628424984894
$ws.Range("E3").Value = -1.155500225866388
$ws.Range("E4").Value = 0.0937874122571985
$ws.Range("E5").Value = 0.04977305232093901
$ws.Range("E6").Value = -0.09861120752912411
$ws.Range("E7").Value = 0.02643600938059715
$ws.Range("E8").Value = -0.05476414598822187
$ws.Range("E9").Value = 0.02129074192242487
$ws.Range("E10").Value = 0.01644016522890406
$ws.Range("E11").Value = 0.01023167260238377
$ws.Range("E12").Value = 0.01689879936321448
